# spreadsheet-data-types.xlsx — "Corrected date handling ... Date Importer tool"
#
# The J column (s="10", custom format "m/d/yyyy;@") holds date-importer
# test fixtures. Row 3's stored serial date (36527) was one day off from
# its siblings in K3/L3 (both 36526) — this is the bug the commit fixes.
# We also leave the worksheet's selection cursor where the author left it
# after verifying the fix (cell J4, just below the corrected value).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the off-by-one date serial in J3 so it matches K3/L3 (36526).
$ws.Range("J3").Value = 36526

# Move/collapse the active selection to J4 (previously J2 with sqref J2:J3).
$ws.Range("J4").Select()
